$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Education section: insert two new bullet points ("CGPA 3.69/4.0 (Cum
#    Laude)" and "Honors: Dean's List (2020 - 2023)") right before the
#    "Electives specializing in Deep Learning and Bioinformatics" bullet,
#    at the same sub-bullet list level.
# ---------------------------------------------------------------------------
$electivesPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "Electives specializing in Deep Learning and Bioinformatics*") {
        $electivesPara = $p
        break
    }
}

$electivesPara.Range.InsertParagraphBefore()
$electivesPara.Range.InsertParagraphBefore()

# Re-find the paragraph so indices line up with the two freshly inserted
# (still empty) paragraphs that now precede it.
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "Electives specializing in Deep Learning and Bioinformatics*") {
        $electivesPara = $p
        break
    }
}

$newPara2 = $electivesPara.Previous()
$newPara1 = $newPara2.Previous()

$newPara1.Range.Text = "CGPA 3.69/4.0 (Cum Laude)"
$newPara2.Range.Text = "Honors: Dean" + [char]8217 + "s List (2020 - 2023)"

# ---------------------------------------------------------------------------
# 2. Research Experience: rework the "digital twin" bullet - demote it to a
#    sub-bullet (ilvl 0 -> 1, indent 720 -> 1440) and replace its text.
# ---------------------------------------------------------------------------
$twinPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*digital twin*") {
        $twinPara = $p
        break
    }
}
$twinPara.Range.ListFormat.ListLevelNumber = 2
$twinPara.LeftIndent = 72

$d.Content.Find.Execute(
    "Studied and developed a digital twin systems (Medical Metaverse) for modeling patient state in real time",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Implemented a digital twin system prototype achieving hundreds of thousands of agents running simultaneously",
    2) | Out-Null

# ---------------------------------------------------------------------------
# 3. Publications: shorten the workshop proceedings name.
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "Proceedings of AI in Oncology Workshop for the International Conference of Artificial Intelligence in Medicine in Europe",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Proceedings of AIME Workshop on AI in Oncology",
    2) | Out-Null

# ---------------------------------------------------------------------------
# 4. Remove the whole "Honors and Awards" section (heading, its separator
#    line, and the three bullet points), from the blank line that follows
#    the TORCS bullet through the "Dean's List Award" bullet. The blank
#    paragraph that precedes "Technical Skills" is left untouched.
# ---------------------------------------------------------------------------
$honorsHeading = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "Honors and Awards*") {
        $honorsHeading = $p
        break
    }
}
$deansListAward = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "Dean*s List Award*") {
        $deansListAward = $p
        break
    }
}

$startPara = $honorsHeading.Previous()
$rangeStart = $startPara.Range.Start
$rangeEnd = $deansListAward.Range.End
$d.Range($rangeStart, $rangeEnd).Delete()
